# "cancel and change booking | purchase seat"
# Adds a new "Change/Cancel" section to the API documentation sheet:
#   row 30 -> section header
#   row 32 -> Cancel booking endpoint (DELETE /booking/id)
#   row 34 -> Change booking endpoint (PUT /booking)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: new section header ---
$ws.Range("A30").Value = "Change/Cancel"

# --- Row 32: Cancel booking API ---
$ws.Rows.Item(32).RowHeight = 61

$ws.Range("A32").Value = "Booking"
$ws.Range("C32").Value = "Cancel"
$ws.Range("D32").Value = "Delete"
$ws.Range("E32").Value = "/booking/id"
$ws.Range("F32").Value = "http://127.0.0.1:5000/booking/619ac7ef63b736cef9773d68"

$cancelSuccess = @"
{
    "message": "Booking #M1ZE6H5CT9PB canceled successfully"
}
"@
$ws.Range("G32").Value = $cancelSuccess

$bookingNotFound = @"
{
    "message": "Booking not found"
}
"@
$ws.Range("H32").Value = $bookingNotFound

$alreadyCanceled = @"
{
    "message": "Booking is already canceled"
}
"@
$ws.Range("I32").Value = $alreadyCanceled

# F32 gets its own distinctive look: Helvetica, size 12, color #212121, wrapped
# (same family/color as the existing Helvetica-8 style used in F10/F11 - copy that
# format across first, then bump the font size to 12)
$ws.Range("F10").Copy()
$ws.Range("F32").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F32").Font.Size = 12
$excel.CutCopyMode = $false

$ws.Range("G32").WrapText = $true
$ws.Range("H32").WrapText = $true
$ws.Range("I32").WrapText = $true

# --- Row 34: Change booking API ---
$ws.Rows.Item(34).RowHeight = 108

$ws.Range("C34").Value = "Change"
$ws.Range("D34").Value = "Put"
$ws.Range("E34").Value = "/booking"

$changePayload = @"
{
    "booking_id":"619ade645bf1fa1382a63771",
    "flight_oid": "619ac7a563b736cef9773d67",
"traveler_details": {"name":"xyz"},
"payment" :{ "reward_points_used" : 20,
  "cash": 50}
}
"@
$ws.Range("F34").Value = $changePayload

$ws.Range("G34").Value = "Booking  changed successfully"

$changeOnce = @"
{
    "message": "You can change the Booking only once!!"
}
"@
$ws.Range("H34").Value = $changeOnce

$ws.Range("F34").WrapText = $true
$ws.Range("G34").WrapText = $true
$ws.Range("H34").WrapText = $true

# --- scroll / selection state, matching where the author ended up editing ---
$ws.Range("F34").Select()
$excel.ActiveWindow.ScrollRow = 23

Write-Output "Change/Cancel booking rows added"
